# Rename the dataset/category value "congenital" to "misc_long_term"
# wherever it appears across all worksheets in the workbook.
#
# Each worksheet in this workbook lists a small set of variable names in
# column A (one of which is "congenital" on some sheets, in either row 2
# or row 3 depending on the sheet). We locate the cell holding the exact
# text "congenital" on every sheet and overwrite it with "misc_long_term".

$wb = $excel.ActiveWorkbook

$updated = 0
foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Cells.Find("congenital", [Type]::Missing, [Type]::Missing, 1)
    if ($cell) {
        $cell.Value = "misc_long_term"
        $updated = $updated + 1
    }
}

Write-Output "Replaced 'congenital' with 'misc_long_term' on $updated worksheet(s)."
